# Refresh the cryptos list (price + 1h volume change) for each coin row.
# Numeric-looking "Price" strings are apostrophe-prefixed so Excel keeps
# storing them as text (matching the original inlineStr cells) instead of
# silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.023.46"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "1.673.75"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'216.16"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "'20.19"
$ws.Range("E10").Value = "  +5.38%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").Value = "1.911.23"
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("D13").Value = "1.663.53"
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "'65.86"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "27.050.69"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "'236.39"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "'7.74"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'4.47"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").Value = "'9.29"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").Value = "'145.76"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'7.15"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "1.470.69"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  +6.28%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "'0.574"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "'0.899"
$ws.Range("E38").Value = "  +7.91%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'6.11"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +12.29%  "
$ws.Range("D43").Value = "'2.27"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("D44").Value = "'66.59"
$ws.Range("E44").Value = "  +7.90%  "
$ws.Range("D45").Value = "1.822.17"
$ws.Range("E45").Value = "  +3.64%  "
$ws.Range("D46").Value = "'0.779"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'90.22"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  +4.23%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "'7.70"
$ws.Range("E51").Value = "  +2.75%  "
